# Update the "dSF" (column F) values to match the re-pulled source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -6
$ws.Range("F5").Value = -6
$ws.Range("F7").Value = -1
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = -5
